$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 5
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 1

$ws.Range("E16").Value = 5
$ws.Range("F16").Value = -5
$ws.Range("G16").Value = 1

$ws.Range("I16").Select()
